$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sprint 4 backlog: row 7 is "View shopping list / Create UI to view current
# shopping list for web application", owned by Destiny. Mark it completed -
# Actual Time (E) = 1, Completed By (F) = Destiny, and the two "Amount
# Remaining After..." week columns (H, I) both drop to 0.
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "Destiny"
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# Move the sheet's active selection to B20 (this also clears the stale
# top-left scroll anchor left over at B1).
$ws.Range("B20").Select()
